# The "usuarios" sheet had a typo in the first user's e-mail address
# (juan.perez@... -> juant.perez@...) and that cell was turned into a
# real mailto: hyperlink. Along with that, the redundant fill formatting
# on the "genero" column data cells was cleared so they share the same
# plain bordered style as the rest of the table, and the selection was
# left on C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the "genero" column (H2:H4) formatting: drop the (invisible,
# but explicitly-flagged) fill so these cells fall back to the common
# bordered/left-aligned style used throughout the sheet.
$ws.Range("H2:H4").Interior.Pattern = -4142   # xlPatternNone

# Fix the e-mail typo for the first user.
$ws.Range("A2").Value = "juant.perez@sena.edu.co"

# Turn that e-mail address into a clickable mailto hyperlink (this also
# applies Excel's built-in "Hyperlink" cell style: underlined, themed
# link color).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:juant.perez@sena.edu.co") | Out-Null

# Leave the selection where the author left it when they saved.
$ws.Range("C6").Select() | Out-Null
